$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "H3" = "0.56"
    "J3" = "0.38"
    "F5" = "0.14"
    "J5" = "0.22"
    "D8" = "0.44"
    "F8" = "0.21"
    "H8" = "0.34"
    "J8" = "0.13"
    "D9" = "0.31"
    "F9" = "0.29"
    "H9" = "0.31"
    "J9" = "0.27"
    "D10" = "0.2"
    "F10" = "0.39"
    "H10" = "0.26"
    "J10" = "0.43"
    "D11" = "0.03"
    "F11" = "0.08"
    "H11" = "0.06"
    "J11" = "0.15"
    "D13" = "3.5"
    "F13" = "2.29"
    "H13" = "3.01"
    "J13" = "2.46"
    "D14" = "0.16"
    "J14" = "0.29"
    "D15" = "0.55"
    "F15" = "0.33"
    "J15" = "0.33"
    "D16" = "0.46"
    "F16" = "0.35"
    "H16" = "0.29"
    "C17" = "621535"
    "D17" = "16159"
    "F17" = "6292"
    "H17" = "1955"
    "J17" = "6260"
    "C19" = "0.55"
    "H19" = "0.6"
    "J19" = "0.45"
    "F20" = "0.28"
    "H20" = "0.28"
    "J20" = "0.35"
    "H21" = "0.11"
    "C22" = "0.02"
    "C24" = "0.57"
    "D24" = "0.43"
    "F24" = "0.18"
    "H24" = "0.27"
    "J24" = "0.11"
    "D25" = "0.31"
    "F25" = "0.31"
    "H25" = "0.3"
    "J25" = "0.25"
    "C26" = "0.06"
    "D26" = "0.19"
    "F26" = "0.38"
    "H26" = "0.31"
    "J26" = "0.43"
    "D27" = "0.05"
    "F27" = "0.09"
    "H27" = "0.1"
    "J27" = "0.19"
    "C29" = "3.71"
    "D29" = "3.62"
    "F29" = "2.43"
    "H29" = "3.18"
    "J29" = "2.5"
    "H30" = "0.19"
    "F31" = "0.3"
    "J31" = "0.27"
    "C32" = "0.75"
    "D32" = "0.75"
    "F32" = "0.62"
    "H32" = "0.68"
    "J32" = "0.68"
    "C33" = "568861"
    "D33" = "13698"
    "F33" = "4681"
    "H33" = "1161"
    "J33" = "4833"
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $changes[$cell]
}
